# Automatische test-sync: 2025-07-31 21:54:50
# Adds a new mail-log entry (row 17) on the "Logs" sheet and bumps the
# matching category counter (B4) on the "Dashboard" sheet from 2 to 3.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 17

$logs.Cells.Item($newRow, 1).Value = "Leg dit even neer bij Koen."
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #15: Leg dit even neer bij Koen."
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$logs.Cells.Item($newRow, 6).Value = "2025-07-31 21:54:44"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(4, 2).Value = 3

# Extend the conditional-formatting ranges so they keep covering the
# data through the newly added row (…2:…16 -> …2:…17) for every
# formatted column.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "16")
    $newRange = $logs.Range($col + "2:" + $col + "17")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
